# "Rerun with new data": the manual traffic count for the 1-hour-interval
# sheet was re-generated. The original rows 8:13 (hours 17-22) now land at
# rows 14:19 (renumbered 20-25), and six brand-new rows (hours 14-19) are
# inserted as the new rows 8:13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Push the existing rows 8:13 down to 14:19, carrying their formatting
#    (styles) along so no new style entries get minted.
$ws.Range("A8:M13").Copy($ws.Range("A14:M19"))

# 2) Renumber column A (the running interval id) for the rows that moved.
$oldIds = @(20, 21, 22, 23, 24, 25)
for ($i = 0; $i -lt $oldIds.Length; $i++) {
    $ws.Cells.Item(14 + $i, 1).Value2 = $oldIds[$i]
}

# 3) Write the freshly (re)counted data into the new rows 8:13.
#    Columns: A  B                   C  D   E  F  G  H   I   J    K  L      M
$newRows = @(
    @(14, 45392.66666666666,  0, 11, 2, 0, 8,  9, 73, 141, 0, 45392, 16),
    @(15, 45392.67361111111,  2,  1, 0, 1, 5,  5, 58, 130, 0, 45392, 16),
    @(16, 45392.68055555555,  2,  4, 1, 1, 3, 14, 56, 179, 0, 45392, 16),
    @(17, 45392.6875,         1,  2, 3, 0, 6,  3, 52, 131, 0, 45392, 16),
    @(18, 45392.69444444445,  1,  7, 1, 0, 6,  7, 54, 165, 0, 45392, 16),
    @(19, 45392.70138888889,  4,  1, 0, 0, 4,  1, 60, 181, 0, 45392, 16)
)

for ($r = 0; $r -lt $newRows.Length; $r++) {
    $row = 8 + $r
    $vals = $newRows[$r]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value2 = $vals[$c]
    }
}

"done"
